# Add four new order sheets (工作表04172024_161026 / _161507 / _161522 / _161554),
# each a duplicate of the existing "出單" snapshot (header + 2 rows) that sheet
# 工作表04172024_155307 already holds, appended at the end of the workbook.

$wb = $excel.ActiveWorkbook

$newSheetNames = @(
    "工作表04172024_161026",
    "工作表04172024_161507",
    "工作表04172024_161522",
    "工作表04172024_161554"
)

foreach ($sheetName in $newSheetNames) {
    $ws = $wb.Worksheets.Add()
    $ws.Name = $sheetName

    # Header row
    $ws.Range("A1").Value = "商品名稱"
    $ws.Range("B1").Value = "下單數量"
    $ws.Range("C1").Value = "售價"
    $ws.Range("D1").Value = "總價格"
    $ws.Range("E1").Value = "負責人"
    $ws.Range("F1").Value = "出單時間"

    # Row 2 - keep "789" in column A as text (matches source data), rest numeric
    $ws.Range("A2").NumberFormat = "@"
    $ws.Range("A2").Value = "789"
    $ws.Range("A2").Style = "Normal"
    $ws.Range("B2").Value = 56
    $ws.Range("C2").Value = 789
    $ws.Range("D2").Value = 44184
    $ws.Range("E2").Value = "AnthonyFu"
    $ws.Range("F2").Value = "2024/04/17 07:53:05"

    # Row 3
    $ws.Range("A3").Value = "傅垣幀"
    $ws.Range("B3").Value = 8
    $ws.Range("C3").Value = 500
    $ws.Range("D3").Value = 4000
    $ws.Range("E3").Value = "AnthonyFu"
    $ws.Range("F3").Value = "2024/04/17 07:53:05"

    # Column widths matching the source sheets (E & F widened to 25)
    $ws.Columns.Item(5).ColumnWidth = 25
    $ws.Columns.Item(6).ColumnWidth = 25

    # Move the freshly created sheet to the end of the workbook (new sheets
    # are inserted at the front by default).
    $ws.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))
}
